# Update "想去人数" (want-to-go count) figures across the workbook's sheets.
# Mirrors the regenerated-data commit: only column F values change.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2523
$ws.Range("F7").Value = 88
$ws.Range("F8").Value = 92
$ws.Range("F11").Value = 1420
$ws.Range("F15").Value = 35
$ws.Range("F16").Value = 995
$ws.Range("F19").Value = 256
$ws.Range("F20").Value = 7611
$ws.Range("F21").Value = 7611
$ws.Range("F22").Value = 8699
$ws.Range("F25").Value = 428
$ws.Range("F31").Value = 11
$ws.Range("F33").Value = 1531
$ws.Range("F39").Value = 40
$ws.Range("F40").Value = 822
$ws.Range("F42").Value = 1379
$ws.Range("F49").Value = 204
$ws.Range("F50").Value = 52

# Sheet "演出" (performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 77
$ws.Range("F19").Value = 325

# Sheet "本地生活" (local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2667
$ws.Range("F5").Value = 164

# Sheet "全部类型" (all types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 164
$ws.Range("F7").Value = 2523
$ws.Range("F9").Value = 92
$ws.Range("F14").Value = 35
$ws.Range("F15").Value = 995
$ws.Range("F19").Value = 77
$ws.Range("F20").Value = 256
$ws.Range("F21").Value = 7611
$ws.Range("F22").Value = 8699
$ws.Range("F24").Value = 428
$ws.Range("F29").Value = 1531
$ws.Range("F35").Value = 40
$ws.Range("F37").Value = 822
$ws.Range("F41").Value = 1379
$ws.Range("F48").Value = 204
$ws.Range("F50").Value = 325
$ws.Range("F51").Value = 52
